$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 434.375
$ws.Range("I2").Value = 123.5
$ws.Range("J2").Value = 745.25
$ws.Range("K2").Value = 123.5
$ws.Range("L2").Value = 745.25
$ws.Range("M2").Value = -10.5
$ws.Range("N2").Value = -971.25
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 318
$ws.Range("I6").Value = 318
$ws.Range("K6").Value = 954
$ws.Range("M6").Value = -842
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 10.5
$ws.Range("I8").Value = 10.5
$ws.Range("K8").Value = 31.5
$ws.Range("M8").Value = 107.5
# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 568.875
$ws.Range("J38").Value = 997.5
$ws.Range("L38").Value = 2992.5
$ws.Range("N38").Value = -3736.5
# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 1632.375
$ws.Range("I58").Value = 296.92307
$ws.Range("J58").Value = 3210.6365
$ws.Range("K58").Value = 890.7692099999999
$ws.Range("L58").Value = 9631.9095
$ws.Range("M58").Value = -740.7692099999999
$ws.Range("N58").Value = -9931.9095
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2863.8
$ws.Range("I132").Value = 1428.5883
$ws.Range("K132").Value = 4285.7649
$ws.Range("M132").Value = -1755.7649

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 6 (Leve Item ID 2226)
$ws.Range("H6").Value = 7001666.5
$ws.Range("I6").Value = 502499.5
$ws.Range("J6").Value = 20000000
$ws.Range("K6").Value = 502499.5
$ws.Range("L6").Value = 20000000
$ws.Range("M6").Value = -502326.5
$ws.Range("N6").Value = -20000346
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2698
$ws.Range("I45").Value = 1831.7778
$ws.Range("K45").Value = 1831.7778
$ws.Range("M45").Value = -1454.7778
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2428.5
$ws.Range("I122").Value = 2428.5
$ws.Range("K122").Value = 7285.5
$ws.Range("M122").Value = -4835.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 2879.6
$ws.Range("I3").Value = 1500
$ws.Range("K3").Value = 1500
$ws.Range("M3").Value = -1386
# Row 5 (Leve Item ID 1750)
$ws.Range("H5").Value = 331.42856
$ws.Range("I5").Value = 303.33334
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 303.33334
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -190.33334
$ws.Range("N5").Value = -726
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 1365
$ws.Range("I7").Value = 797.5
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 797.5
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -684.5
$ws.Range("N7").Value = -2726
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3363.111
$ws.Range("I31").Value = 2829.3333
$ws.Range("J31").Value = 4430.6665
$ws.Range("K31").Value = 2829.3333
$ws.Range("L31").Value = 4430.6665
$ws.Range("M31").Value = -2534.3333
$ws.Range("N31").Value = -5020.6665
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3363.111
$ws.Range("I34").Value = 2829.3333
$ws.Range("J34").Value = 4430.6665
$ws.Range("K34").Value = 2829.3333
$ws.Range("L34").Value = 4430.6665
$ws.Range("M34").Value = -2627.3333
$ws.Range("N34").Value = -4834.6665
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 3019.7334
$ws.Range("I58").Value = 1941
$ws.Range("J58").Value = 7334.6665
$ws.Range("K58").Value = 1941
$ws.Range("L58").Value = 7334.6665
$ws.Range("M58").Value = -1738
$ws.Range("N58").Value = -7740.6665
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 13140.286
$ws.Range("I99").Value = 7766.385
$ws.Range("J99").Value = 17797.666
$ws.Range("K99").Value = 7766.385
$ws.Range("L99").Value = 17797.666
$ws.Range("M99").Value = -6268.385
$ws.Range("N99").Value = -20793.666
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 2987.3333
$ws.Range("I122").Value = 2987.3333
$ws.Range("K122").Value = 8961.999899999999
$ws.Range("M122").Value = -6511.999899999999
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 13140.286
$ws.Range("I126").Value = 7766.385
$ws.Range("J126").Value = 17797.666
$ws.Range("K126").Value = 23299.155
$ws.Range("L126").Value = 53392.99800000001
$ws.Range("M126").Value = -20829.155
$ws.Range("N126").Value = -58332.99800000001
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 3019.7334
$ws.Range("I136").Value = 1941
$ws.Range("J136").Value = 7334.6665
$ws.Range("K136").Value = 5823
$ws.Range("L136").Value = 22003.9995
$ws.Range("M136").Value = -3273
$ws.Range("N136").Value = -27103.9995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 285.26666
$ws.Range("I2").Value = 46.25
$ws.Range("J2").Value = 372.18182
$ws.Range("K2").Value = 277.5
$ws.Range("L2").Value = 2233.09092
$ws.Range("M2").Value = -164.5
$ws.Range("N2").Value = -2459.09092
# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 33333446
$ws.Range("I7").Value = 33333446
$ws.Range("K7").Value = 100000338
$ws.Range("M7").Value = -100000226
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 10611.214
$ws.Range("J39").Value = 15632.125
$ws.Range("L39").Value = 46896.375
$ws.Range("N39").Value = -47484.375
# Row 104 (Leve Item ID 19807)
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 125 (Leve Item ID 36043)
$ws.Range("H125").Value = 3000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 9000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -18840
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 1128.75
$ws.Range("I129").Value = 776
$ws.Range("J129").Value = 1716.6666
$ws.Range("K129").Value = 2328
$ws.Range("L129").Value = 5149.9998
$ws.Range("M129").Value = 2672
$ws.Range("N129").Value = -15149.9998
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1397.25
$ws.Range("J131").Value = 1385.6875
$ws.Range("L131").Value = 4157.0625
$ws.Range("N131").Value = -14237.0625
# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 1666.3334
$ws.Range("I134").Value = 1666.3334
$ws.Range("K134").Value = 4999.0002
$ws.Range("M134").Value = 70.9997999999996
# Row 139 (Leve Item ID 44102)
$ws.Range("H139").Value = 2771
$ws.Range("I139").Value = 1417.6666
$ws.Range("J139").Value = 5477.6665
$ws.Range("K139").Value = 4252.9998
$ws.Range("L139").Value = 16432.9995
$ws.Range("M139").Value = 887.0002000000004
$ws.Range("N139").Value = -26712.9995

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 22 (Leve Item ID 2685)
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1966.6666
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 1966.6666
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -1437.6666
$ws.Range("N22").Value = -1158
# Row 118 (Leve Item ID 26172)
$ws.Range("H118").Value = 43000
$ws.Range("J118").Value = 43000
$ws.Range("L118").Value = 43000
$ws.Range("N118").Value = -46314
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 5763
$ws.Range("I126").Value = 5581.6665
$ws.Range("J126").Value = 6035
$ws.Range("K126").Value = 16744.9995
$ws.Range("L126").Value = 18105
$ws.Range("M126").Value = -14274.9995
$ws.Range("N126").Value = -23045
# Row 129 (Leve Item ID 35367)
$ws.Range("H129").Value = 21999.5
$ws.Range("J129").Value = 21999.5
$ws.Range("L129").Value = 21999.5
$ws.Range("N129").Value = -31999.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2462.125
$ws.Range("I100").Value = 2533
$ws.Range("K100").Value = 2533
$ws.Range("M100").Value = -1992

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 117 (Leve Item ID 26162)
$ws.Range("H117").Value = 30409
$ws.Range("J117").Value = 30409
$ws.Range("L117").Value = 30409
$ws.Range("N117").Value = -39587
